# feat: add 2022-Q1 data
#
# The previous "总计" sheet (sheetId 3 / rId3) is renamed to "2022-Q1" and
# repopulated with the new quarter's per-fund holding detail (same shape
# as the existing "2021-Q4" sheet). A brand-new "总计" sheet (sheetId 4 /
# rId4) is appended right after it, holding the refreshed
# quarter-over-quarter summary table (new "2022-Q1" row on top, the
# previous two rows pushed down).

$wb = $excel.ActiveWorkbook

# Template sheet whose header-row / index-column formatting (style index
# "2": bold, centered, thin-bordered) we reuse for the two sheets below.
$fmtSrc = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1. Turn the previous "总计" sheet into the new "2022-Q1" detail sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")

# Wipe the old 3-row "总计" content (A1:D3) before laying out the new table.
$q1.Range("A1:D3").Clear()

$q1.Name = "2022-Q1"

# Header row.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Make sure the fund-code / numeric-look text columns are stored as plain
# text (so leading zeros like "005443" survive), same as the other
# per-quarter sheets. Flip to a text format just long enough to type the
# values in, then drop back to General so no stray per-cell style is
# left behind (matches the source, which has no "s" attribute on these
# data cells).
$q1BodyRange = $q1.Range("B2:G5")
$q1BodyRange.NumberFormat = "@"

$q1Data = @(
    @(0, "519613", "银河君尚灵活配置混合A",       "5.70", "29.99", "0.79", "0.0450", 2),
    @(1, "519615", "银河君尚灵活配置混合I",       "5.70", "29.99", "0.79", "0.0450", 2),
    @(2, "005443", "国金量化多策略灵活配置混合", "0.51", "64.10", "0.90", "0.0046", 5),
    @(3, "519614", "银河君尚灵活配置混合C",       "0.23", "29.99", "0.79", "0.0018", 2)
)

$r = 2
foreach ($row in $q1Data) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = $row[3]
    $q1.Cells.Item($r, 5).Value = $row[4]
    $q1.Cells.Item($r, 6).Value = $row[5]
    $q1.Cells.Item($r, 7).Value = $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

$q1BodyRange.Style = "Normal"

# Re-apply the bold/centered/bordered header + index-column look (lost
# when we cleared the old "总计" range above) by copying formats from the
# equivalent cells on "2021-Q4".
$fmtSrc.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$fmtSrc.Range("A2").Copy()
$q1.Range("A2:A5").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Append a fresh "总计" sheet right after "2022-Q1" with the updated
#    quarter summary (new quarter on top, older ones pushed down).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# Match the page-margin layout the other sheets (and the previous "总计"
# sheet) use, rather than the Add() default.
$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalData = @(
    @(0, "2022-Q1", 4, 0.1),
    @(1, "2021-Q4", 3, 0.1),
    @(2, "2021-Q3", 5, 0.07000000000000001)
)

$r = 2
foreach ($row in $totalData) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$fmtSrc.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$fmtSrc.Range("A2").Copy()
$total.Range("A2:A4").PasteSpecial(-4122)
